$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "44.086.97"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.63%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.265.17"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.31%  "

# Row 4
$ws.Range("E4").Value = "  +0.12%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "231.03"
$ws.Range("D5").Style = "Normal"

# Row 6
$ws.Range("E6").Value = "  +1.74%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "64.59"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.86%  "

# Row 8
$ws.Range("E8").Value = "  +0.00%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.449"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.95%  "

# Row 10
$ws.Range("E10").Value = "  +5.26%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "57.24"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.93%  "

# Row 12
$ws.Range("E12").Value = "  +16.34%  "

# Row 13
$ws.Range("E13").Value = "  +1.95%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.604.04"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.14%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.77"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.12%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.11"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.06%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.840"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.64%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.284.42"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.77%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "44.006.97"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.74%  "

# Row 20
$ws.Range("E20").Value = "  +7.37%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "73.84"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.06%  "

# Row 22
$ws.Range("E22").Value = "  -1.88%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "253.35"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.50%  "

# Row 24
$ws.Range("E24").Value = "  +0.11%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.45"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.62%  "

# Row 26
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.31"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.52%  "

# Row 27
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.12"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.75%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.29"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +23.48%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "171.68"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.52%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.141"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.07%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.94"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.68%  "

# Row 32
$ws.Range("E32").Value = "  -4.04%  "

# Row 33
$ws.Range("E33").Value = "  +3.04%  "

# Row 34
$ws.Range("E34").Value = "  +6.56%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.82"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.64%  "

# Row 36
$ws.Range("E36").Value = "  -3.55%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.81"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.73%  "

# Row 38
$ws.Range("E38").Value = "  +0.54%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.33"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.52%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0261"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.95%  "

# Row 41
$ws.Range("E41").Value = "  +0.05%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.000225"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.46%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0988"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.11%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "17.54"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.08%  "

# Row 45
$ws.Range("B45").Value = "Celestia"
$ws.Range("C45").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.58"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +8.23%  "

# Row 46
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.24"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.27%  "

# Row 47
$ws.Range("B47").Value = "TrustWalletToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.21"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.55%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "98.47"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.40%  "

# Row 49
$ws.Range("E49").Value = "  -2.16%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.37"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.25%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.448.75"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.66%  "
